$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("System")

$ws.Range("B3").Value = 50
$ws.Range("C3").Value = 10
